$d = $word.ActiveDocument

# 1. Update activation date: 01/01/2016 -> 01/01/2023
$d.Content.Find.Execute("Ativação: 01/01/2016", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2023", 2)

# English translations to be added as new italic paragraphs right after their
# Portuguese counterparts.
$objetivosEn = "Develop theoretical and practical knowledge of the manufacturing processes of equipment and devices required for the development of products and prototypes. Know the requirements and effects of manufacturing processes in order to allow, interact, create and execute projects throughout your professional life."
$programaResumidoEn = "Introduction to manufacturing processes. Material joining processes. Computer-aided design (CAD) review. Computer Aided Manufacturing (CAM). Flexible production systems. Rapid prototyping."
$programaEn = "Classification of manufacturing processes. Foundry. Powder metallurgy. Machining: processes, fundamentals and economic conditions. Machine tools. Mechanical conformation. Material joining processes. Computer Aided Manufacturing (CAM). Programming languages for numerical control. Numerical control machine tools. Product manufacturing sequence. Notions of automation of manufacturing processes. Rapid prototyping. Rapid prototyping systems (solid, liquid and powder)."

function Insert-ItalicParagraphAfter($matchPrefix, $englishText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($matchPrefix)) {
            $p.Range.InsertParagraphAfter()
            $newRange = $d.Paragraphs.Item($i + 1).Range
            # exclude the paragraph mark itself so only the run gets italics
            $newRange.MoveEnd(1, -1)
            $newRange.Text = $englishText
            $newRange.Font.Italic = $true
            return
        }
    }
}

# Insert from the bottom of the document upward so that paragraph indices for
# targets not yet processed remain unaffected by earlier insertions.
Insert-ItalicParagraphAfter "Classificação dos processos de fabricação" $programaEn
Insert-ItalicParagraphAfter "Introdução aos processos de fabricação" $programaResumidoEn
Insert-ItalicParagraphAfter "Desenvolver conhecimento teórico" $objetivosEn
